$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet Sheet1 -> Address
$ws.Name = "Address"

# 2. Fix H2/H3/H4 mobile numbers: re-enter as text (quote-prefixed) so the
#    leading zeros / large numeric strings round-trip like real phone numbers
$ws.Range("H2").Value = "'9899001068"
$ws.Range("H3").Value = "'9899001068"
$ws.Range("H4").Value = "'1234567890"

# 3. Fill row 5 with the new address record (previously blank placeholder row)
$ws.Range("A5").Value = "Rohan Shukla"
$ws.Range("A5").Interior.ColorIndex = -4142

$ws.Range("B5").Value = "Mumbai"
$ws.Range("B5").Interior.ColorIndex = -4142

$ws.Range("C5").Value = "mumbai2"
$ws.Range("C5").Interior.ColorIndex = -4142

$ws.Range("D5").Value = "Seepz"
$ws.Range("D5").Interior.ColorIndex = -4142

$ws.Range("E5").Value = "Andheri"
$ws.Range("E5").Interior.ColorIndex = -4142

$ws.Range("F5").Value = "Maharashtra"

$ws.Range("G5").Value = 400096
$ws.Range("G5").Interior.ColorIndex = -4142

$ws.Range("H5").Value = "'9876543210"
$ws.Range("H5").Interior.ColorIndex = -4142

# 4. Drop the now-unused trailing blank row 6
$ws.Range("A6:H6").Clear()

# 5. Match source workbook's saved selection
$ws.Range("E11").Select()
